$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for every
# coin row, plus a handful of coin re-rankings where the name/link/price/volume
# of two rows trade places (rows 33/34, 47/48) or a coin is fully replaced by a
# new entrant (row 51: Stacks -> WEMIXToken).
#
# Price values are prefixed with a leading apostrophe so Excel keeps them as
# literal text (preserving formats like "1.00", "579.60", "66.442.22") instead
# of auto-converting them to floating point numbers.

$ws.Range("D2").Value = "'66.442.22"
$ws.Range("E2").Value = "  -3.50%  "

$ws.Range("D3").Value = "'3.562.31"
$ws.Range("E3").Value = "  -3.75%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'579.60"
$ws.Range("E5").Value = "  -5.67%  "

$ws.Range("D6").Value = "'184.65"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("D7").Value = "'3.559.03"
$ws.Range("E7").Value = "  -3.68%  "

$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.667"
$ws.Range("E10").Value = "  -6.07%  "

$ws.Range("D11").Value = "'0.145"
$ws.Range("E11").Value = "  -8.42%  "

$ws.Range("D12").Value = "'52.66"
$ws.Range("E12").Value = "  -4.92%  "

$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  -9.77%  "

$ws.Range("D14").Value = "'9.76"
$ws.Range("E14").Value = "  -6.74%  "

$ws.Range("D15").Value = "'4.125.43"
$ws.Range("E15").Value = "  -3.76%  "

$ws.Range("D16").Value = "'3.566.36"
$ws.Range("E16").Value = "  -3.62%  "

$ws.Range("E17").Value = "  -0.80%  "

$ws.Range("D18").Value = "'18.26"
$ws.Range("E18").Value = "  -4.75%  "

$ws.Range("D19").Value = "'12.14"
$ws.Range("E19").Value = "  -5.47%  "

$ws.Range("D20").Value = "'66.321.13"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("E21").Value = "  -6.04%  "

$ws.Range("D22").Value = "'393.79"
$ws.Range("E22").Value = "  -3.58%  "

$ws.Range("D23").Value = "'4.31"
$ws.Range("E23").Value = "  -5.73%  "

$ws.Range("D24").Value = "'85.61"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("D25").Value = "'11.11"
$ws.Range("E25").Value = "  +2.18%  "

$ws.Range("D26").Value = "'2.89"
$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").Value = "'12.37"
$ws.Range("E27").Value = "  -2.58%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'3.52"
$ws.Range("E29").Value = "  -5.79%  "

$ws.Range("D30").Value = "'8.89"
$ws.Range("E30").Value = "  -7.08%  "

$ws.Range("D31").Value = "'30.98"
$ws.Range("E31").Value = "  -5.71%  "

$ws.Range("D32").Value = "'7.04"
$ws.Range("E32").Value = "  -1.69%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'626.60"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'12.15"
$ws.Range("E34").Value = "  -3.00%  "

$ws.Range("D35").Value = "'63.79"
$ws.Range("E35").Value = "  -2.43%  "

$ws.Range("E36").Value = "  -6.66%  "

$ws.Range("D37").Value = "'41.10"
$ws.Range("E37").Value = "  -6.41%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'0.397"
$ws.Range("E39").Value = "  -3.03%  "

$ws.Range("D40").Value = "'0.0₃0774"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("E41").Value = "  -6.01%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'2.972.47"
$ws.Range("E43").Value = "  +4.12%  "

$ws.Range("E44").Value = "  -6.49%  "

$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").Value = "'0.0406"
$ws.Range("E46").Value = "  -7.54%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.130"
$ws.Range("E47").Value = "  -6.29%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.09"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").Value = "'8.45"
$ws.Range("E49").Value = "  -6.43%  "

$ws.Range("D50").Value = "'137.38"
$ws.Range("E50").Value = "  -3.05%  "

$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.50"
$ws.Range("E51").Value = "  -8.29%  "
